$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.832.12"
$ws.Range("E2").Value = "  -2.69%  "
$ws.Range("D3").Value = "3.203.75"
$ws.Range("E3").Value = "  -4.50%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "531.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.90%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "3.201.17"
$ws.Range("E9").Value = "  -4.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.608"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -8.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.133"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000252"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").Value = "3.714.61"
$ws.Range("E15").Value = "  -4.06%  "
$ws.Range("E16").Value = "  -3.47%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "17.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.187.56"
$ws.Range("E18").Value = "  -4.75%  "
$ws.Range("D19").Value = "62.693.76"
$ws.Range("E19").Value = "  -2.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.968"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "367.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "643.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.34"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.106"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.17%  "
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.377"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.995"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  +14.28%  "
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("D43").Value = "2.877.50"
$ws.Range("E43").Value = "  +3.10%  "
$ws.Range("E44").Value = "  +9.14%  "
$ws.Range("E45").Value = "  +10.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.49%  "
$ws.Range("E47").Value = "  +1.22%  "
$ws.Range("E48").Value = "  -4.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.48%  "
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.54%  "
